$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the MobileNumberField_Xpath value (row 13, column B)
$ws.Range("B13").Value = "//input[@name='socialogin_email']"

# Restyle rows 12-15 (existing cells) to match the rest of the Courier-New /
# Calibri scheme used throughout the sheet.
# A12:A15 -> Calibri 11 black (matches A4:A11)
$ws.Range("A12:A15").Font.Name = "Calibri"
$ws.Range("A12:A15").Font.Size = 11
$ws.Range("A12:A15").Font.Color = 0

# B12 -> Calibri 11 black (matches B7/B11 "header-like" cells)
$ws.Range("B12").Font.Name = "Calibri"
$ws.Range("B12").Font.Size = 11
$ws.Range("B12").Font.Color = 0

# B13:B15 -> Courier New 10 blue (matches B2,B4..B10 xpath value cells)
$ws.Range("B13:B15").Font.Name = "Courier New"
$ws.Range("B13:B15").Font.Size = 10
$ws.Range("B13:B15").Font.Color = 16711722

# Add new rows 16-24
$data = @(
    @("Login_For_Review_And_Rating_LinkText", "Log in"),
    @("Rating_Text_Xpath", "(//form[@id='review-form']//div[2])[1]/h5"),
    @("Review_Text_Xpath", "(//form[@id='review-form']//div[2])[1]/ul/h5"),
    @("Rating_Stars_Xpath", "(//form[@id='review-form']//div[2])[1]/div[1]/ul/li"),
    @("Empty_Rating_Stars_Xpath", "(//form[@id='review-form']//div[2])[1]/div[1]/div[2]/input"),
    @("Review_Title_Label_Xpath", "//*[@id='review-form']/fieldset/div[2]/ul/li[1]/label"),
    @("Review_Title_TextBox_Xpath", "//*[@id='summary_field']"),
    @("Review_Xpath", "//*[@id='review-form']/fieldset/div[2]/ul/li[2]/label"),
    @("Review_TextBox_Xpath", "//*[@id='review_field']")
)

$row = 16
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# New rows use the plain Calibri look, except A16 which uses the Courier-New
# blue look (matches the "key name is itself an xpath-ish label" pattern
# seen at A13/A24 etc. in the target).
$ws.Range("A16:A24").Font.Name = "Calibri"
$ws.Range("A16:A24").Font.Size = 11
$ws.Range("A16:A24").Font.Color = 0
$ws.Range("B16:B24").Font.Name = "Calibri"
$ws.Range("B16:B24").Font.Size = 11
$ws.Range("B16:B24").Font.Color = 0

$ws.Range("A16").Font.Name = "Courier New"
$ws.Range("A16").Font.Size = 10
$ws.Range("A16").Font.Color = 16711722

$ws.Range("B25").Select()
